$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.120.01"
$ws.Range('E2').Value = '  -0.49%  '

$ws.Range('D3').Value = "'1.813.45"
$ws.Range('E3').Value = '  +1.65%  '

$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  -0.43%  '

$ws.Range('D5').Value = "'337.89"
$ws.Range('E5').Value = '  -0.74%  '

$ws.Range('D6').Value = "'0.9990"
$ws.Range('E6').Value = '  -0.25%  '

$ws.Range('D7').Value = "'0.4021"
$ws.Range('E7').Value = '  +5.05%  '

$ws.Range('D8').Value = "'0.3456"
$ws.Range('E8').Value = '  +0.53%  '

$ws.Range('D9').Value = "'45.93"
$ws.Range('E9').Value = '  -2.53%  '

$ws.Range('D10').Value = "'1.157"
$ws.Range('E10').Value = '  +0.15%  '

$ws.Range('D11').Value = "'0.07436"
$ws.Range('E11').Value = '  +0.41%  '

$ws.Range('D12').Value = "'22.89"
$ws.Range('E12').Value = '  -3.47%  '

$ws.Range('D13').Value = "'0.9999"
$ws.Range('E13').Value = '  -0.21%  '

$ws.Range('D14').Value = "'6.279"
$ws.Range('E14').Value = '  -2.89%  '

$ws.Range('D15').Value = "'7.309"
$ws.Range('E15').Value = '  -0.92%  '

$ws.Range('D16').Value = "'1.807.52"
$ws.Range('E16').Value = '  +0.85%  '

$ws.Range('D17').Value = "'0.00001082"
$ws.Range('E17').Value = '  +0.47%  '

$ws.Range('D18').Value = "'0.06660"
$ws.Range('E18').Value = '  -0.43%  '

$ws.Range('D19').Value = "'82.16"
$ws.Range('E19').Value = '  -0.16%  '

$ws.Range('D20').Value = "'0.9977"
$ws.Range('E20').Value = '  -0.39%  '

$ws.Range('E21').Value = '  -1.30%  '

$ws.Range('D22').Value = "'6.368"
$ws.Range('E22').Value = '  -0.72%  '

$ws.Range('D23').Value = "'28.150.30"
$ws.Range('E23').Value = '  -0.49%  '

$ws.Range('D24').Value = "'11.88"
$ws.Range('E24').Value = '  -2.02%  '

$ws.Range('D25').Value = "'2.406"
$ws.Range('E25').Value = '  +1.34%  '

$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = "'2.450"
$ws.Range('E26').Value = '  +1.56%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'20.67"
$ws.Range('E27').Value = '  -0.65%  '

$ws.Range('D28').Value = "'155.28"
$ws.Range('E28').Value = '  +0.48%  '

$ws.Range('D29').Value = "'2.012.71"
$ws.Range('E29').Value = '  +0.23%  '

$ws.Range('D30').Value = "'1.324"
$ws.Range('E30').Value = '  -7.09%  '

$ws.Range('D31').Value = "'132.16"
$ws.Range('E31').Value = '  -2.39%  '

$ws.Range('D32').Value = "'4.076"
$ws.Range('E32').Value = '  +1.39%  '

$ws.Range('D33').Value = "'6.003"
$ws.Range('E33').Value = '  -1.84%  '

$ws.Range('D34').Value = "'0.08822"
$ws.Range('E34').Value = '  -1.35%  '

$ws.Range('D35').Value = "'12.40"
$ws.Range('E35').Value = '  -2.97%  '

$ws.Range('D36').Value = "'0.02344"
$ws.Range('E36').Value = '  -2.90%  '

$ws.Range('D37').Value = "'0.06317"
$ws.Range('E37').Value = '  -1.02%  '

$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D38').Value = "'0.6672"
$ws.Range('E38').Value = '  -2.77%  '

$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = "'5.248"
$ws.Range('E39').Value = '  -2.20%  '

$ws.Range('E40').Value = '  -0.53%  '

$ws.Range('D41').Value = "'1.520"
$ws.Range('E41').Value = '  +1.17%  '

$ws.Range('E42').Value = '  -2.65%  '

$ws.Range('D43').Value = "'8.144"
$ws.Range('E43').Value = '  -2.08%  '

$ws.Range('D44').Value = "'14.21"
$ws.Range('E44').Value = '  -0.01%  '

$ws.Range('D45').Value = "'0.9980"
$ws.Range('E45').Value = '  -0.33%  '

$ws.Range('D46').Value = "'0.6160"
$ws.Range('E46').Value = '  -2.12%  '

$ws.Range('D47').Value = "'3.877"
$ws.Range('E47').Value = '  -0.20%  '

$ws.Range('D48').Value = "'128.54"
$ws.Range('E48').Value = '  -3.52%  '

$ws.Range('D49').Value = "'2.052"
$ws.Range('E49').Value = '  -1.48%  '

$ws.Range('D50').Value = "'1.177"
$ws.Range('E50').Value = '  -1.57%  '

$ws.Range('D51').Value = "'0.07124"
$ws.Range('E51').Value = '  -4.66%  '
